# Full Project including advanced topic
#
# - Sheet1 ("login" sheet): admin/admin123 -> tutorial/tutorial row1,
#   Admin/admin123 row2; selection moves to A2.
# - Sheet2: turned from an empty sheet into a one-row FMLA record
#   (name, type, month/year/day start, month/year/day end, description),
#   autofit columns, portrait page orientation, tab made active,
#   selection on H1.
# - Sheet3: removed entirely.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 updates first (reuses existing shared strings "admin123"/
#     "tutorial"; orphans the old "admin" string so it is dropped when
#     the shared-string table is compacted on save, and interns "Admin"
#     as the next new string).
$ws1.Range("A1").Value = "tutorial"
$ws1.Range("B1").Value = "tutorial"
$ws1.Range("A2").Value = "Admin"
$ws1.Range("B2").Value = "admin123"

# --- Sheet2 content next, in the same order the target workbook
#     introduced each new string (Dec, This is a test, Cassidy Hope,
#     US - FMLA).
$ws2.Range("C1").Value = "Dec"
$ws2.Range("D1").Value = 2020
$ws2.Range("E1").Value = 14
$ws2.Range("F1").Value = "Dec"
$ws2.Range("G1").Value = 2020
$ws2.Range("H1").Value = 19
$ws2.Range("I1").Value = "This is a test"
$ws2.Range("A1").Value = "Cassidy Hope"
$ws2.Range("B1").Value = "US - FMLA"

# --- Remove Sheet3
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Delete()

# --- Column widths / page layout for Sheet2
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(7).AutoFit()
$ws2.PageSetup.Orientation = 1

# --- Selections / active sheet
$ws1.Range("A2").Select()
$ws2.Range("H1").Select()
$ws2.Activate()
